$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: take on the data previously held by row 4
$ws.Range("D2").Value = 44355
$ws.Range("L2").Value = "Segunda"
$ws.Range("N2").Value = 20000
$ws.Range("O2").Value = 21000
$ws.Range("P2").Value = 20500
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 1139

# Row 4: take on the data previously held by row 5
$ws.Range("D4").Value = 44342
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 250
$ws.Range("N4").Value = 24000
$ws.Range("O4").Value = 25000
$ws.Range("P4").Value = 24500
$ws.Range("R4").Value = "Región de O'Higgins"
$ws.Range("S4").Value = 1361

# Row 5: take on the data previously held by row 2
$ws.Range("D5").Value = 44313
$ws.Range("M5").Value = 270
$ws.Range("N5").Value = 21000
$ws.Range("O5").Value = 22000
$ws.Range("P5").Value = 21500
$ws.Range("S5").Value = 1194
